$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.753.50"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.246.42"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0832"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").Value = "2.588.78"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.855"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.291.94"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "43.645.09"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.18%  "
$ws.Range("D20").Value = "0.0₃0987"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("E35").Value = "  +8.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("E37").Value = "  -3.46%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +22.74%  "
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "1.808.72"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "75.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
